$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 338, shifting existing rows 338:353 down to 339:354.
$ws.Rows(338).Insert()

# Populate the newly inserted row 338 with a new weekly price record
# (same Mercado/Categoria/etc. as its neighboring rows, new Fecha + prices).
$ws.Range("A338").Value = 8
$ws.Range("B338").Value = "Terminal La Palmera de La Serena"
$ws.Range("C338").Value = "Coquimbo"
$ws.Range("D338").Value = 45147
$ws.Range("E338").Value = 4
$ws.Range("F338").Value = 100112037
$ws.Range("G338").Value = "Cebollín"
$ws.Range("H338").Value = "Sin especificar"
$ws.Range("I338").Value = "Primera"
$ws.Range("J338").Value = 1000
$ws.Range("K338").Value = 1000
$ws.Range("L338").Value = 1200
$ws.Range("M338").Value = 1100
$ws.Range("N338").Value = "$/paquete 6 unidades"
$ws.Range("O338").Value = "Provincia del Elquí"
$ws.Range("P338").Value = 183
$ws.Range("Q338").Value = 6
$ws.Range("R338").Value = "Hortaliza"
